# "Card List" workbook update:
#   - Insert a blank separator row before each change of card "Type"
#     (Minion -> Spell -> Enchantment -> Ritual) in the card table.
#   - Update the active selection.
#   - Set the page orientation to portrait.
#
# Working from the bottom of the sheet upward so each insertion uses the
# original (pre-shift) row numbers and doesn't disturb rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank row before "Dark Ritual" (start of the Ritual section).
$ws.Rows.Item(21).Insert()

# Blank row before "Giant Strength" (start of the Enchantment section).
$ws.Rows.Item(16).Insert()

# Blank row before "Banish" (start of the Spell section).
$ws.Rows.Item(10).Insert()

# Update the selected cell to reflect where the user left off.
[void]$ws.Range("J27").Select()

# Set the page to portrait orientation.
$ws.PageSetup.Orientation = 1
